# Add "Norway" and "Poland" market test-data sheets after "Hungary",
# cloned from the existing "Turkey" sheet template (same layout/column
# widths as the other market sheets), then fill in the per-country values.

$wb = $excel.ActiveWorkbook

$turkey  = $wb.Worksheets.Item("Turkey")
$hungary = $wb.Worksheets.Item("Hungary")

# --- Norway: copy Turkey's sheet right after Hungary, rename, fill values ---
$turkey.Copy($null, $hungary)
$norway = $wb.Worksheets.Item($hungary.Index + 1)
$norway.Name = "Norway"
$norway.Range("B4").Value = "NGC-2931/T3073/T3071"
$norway.Range("B2").Value = "Norway Market"

# --- Poland: copy Turkey's sheet right after Norway, rename, fill values ---
$turkey.Copy($null, $norway)
$poland = $wb.Worksheets.Item($norway.Index + 1)
$poland.Name = "Poland"
$poland.Range("B4").Value = "NGC-2920/T3039/T3037"
$poland.Range("B2").Value = "Poland Market"

# Norway ends up the active/selected tab.
$norway.Activate()
